# Auto-generated Excel COM-interop script
# Applies the meteocat daily-summary refresh diff (2026-02-17 03:51 run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Humitat (%) cells must be forced to Text format first, otherwise Excel's
# auto-detection parses '48%' etc. as a numeric percentage instead of literal text.
$percentCells = @("H2", "H6", "H7", "H8", "H9", "H13", "H14", "H16", "H17", "H19", "H20", "H23", "H28", "H29", "H30", "H34", "H36", "H37", "H39", "H40", "H44", "H45", "H46")
foreach ($addr in $percentCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("E2").Value = '2026-02-17 03:49:04'
$ws.Range("H2").Value = '48%'
$ws.Range("E3").Value = '2026-02-17 03:49:06'
$ws.Range("N3").Value = '-8.2 °C 3:02 TU'
$ws.Range("O3").Value = '-7.0 °C'
$ws.Range("E4").Value = '2026-02-17 03:49:09'
$ws.Range("J4").Value = '1014.9 hPa'
$ws.Range("E5").Value = '2026-02-17 03:49:12'
$ws.Range("O5").Value = '-6.5 °C'
$ws.Range("E6").Value = '2026-02-17 03:49:15'
$ws.Range("H6").Value = '84%'
$ws.Range("J6").Value = '1014.6 hPa'
$ws.Range("O6").Value = '8.6 °C'
$ws.Range("E7").Value = '2026-02-17 03:49:17'
$ws.Range("H7").Value = '53%'
$ws.Range("L7").Value = '38.5 km/h - 322º 3:09 TU'
$ws.Range("N7").Value = '13.8 °C 3:29 TU'
$ws.Range("O7").Value = '14.7 °C'
$ws.Range("E8").Value = '2026-02-17 03:49:20'
$ws.Range("H8").Value = '70%'
$ws.Range("J8").Value = '1014.4 hPa'
$ws.Range("N8").Value = '9.5 °C 3:29 TU'
$ws.Range("O8").Value = '10.4 °C'
$ws.Range("E9").Value = '2026-02-17 03:49:22'
$ws.Range("H9").Value = '44%'
$ws.Range("N9").Value = '9.6 °C 3:06 TU'
$ws.Range("O9").Value = '12.7 °C'
$ws.Range("E10").Value = '2026-02-17 03:49:25'
$ws.Range("N10").Value = '6.8 °C 3:19 TU'
$ws.Range("O10").Value = '7.6 °C'
$ws.Range("E11").Value = '2026-02-17 03:49:28'
$ws.Range("N11").Value = '5.2 °C 3:29 TU'
$ws.Range("O11").Value = '8.1 °C'
$ws.Range("E12").Value = '2026-02-17 03:49:30'
$ws.Range("O12").Value = '13.1 °C'
$ws.Range("E13").Value = '2026-02-17 03:49:33'
$ws.Range("H13").Value = '41%'
$ws.Range("J13").Value = '1016.0 hPa'
$ws.Range("N13").Value = '1.8 °C 3:10 TU'
$ws.Range("O13").Value = '4.9 °C'
$ws.Range("E14").Value = '2026-02-17 03:49:36'
$ws.Range("H14").Value = '60%'
$ws.Range("L14").Value = '35.3 km/h - 286º 3:07 TU'
$ws.Range("N14").Value = '12.7 °C 3:00 TU'
$ws.Range("O14").Value = '13.6 °C'
$ws.Range("E15").Value = '2026-02-17 03:49:38'
$ws.Range("N15").Value = '11.9 °C 3:29 TU'
$ws.Range("O15").Value = '12.6 °C'
$ws.Range("E16").Value = '2026-02-17 03:49:41'
$ws.Range("H16").Value = '46%'
$ws.Range("E17").Value = '2026-02-17 03:49:44'
$ws.Range("G17").Value = '1 cm'
$ws.Range("H17").Value = '54%'
$ws.Range("N17").Value = '1.3 °C 3:20 TU'
$ws.Range("O17").Value = '2.7 °C'
$ws.Range("E18").Value = '2026-02-17 03:49:46'
$ws.Range("J18").Value = '1014.9 hPa'
$ws.Range("N18").Value = '5.6 °C 3:15 TU'
$ws.Range("E19").Value = '2026-02-17 03:49:49'
$ws.Range("H19").Value = '72%'
$ws.Range("L19").Value = '21.2 km/h - 11º 3:02 TU'
$ws.Range("N19").Value = '5.4 °C 3:28 TU'
$ws.Range("E20").Value = '2026-02-17 03:49:51'
$ws.Range("H20").Value = '44%'
$ws.Range("N20").Value = '-4.1 °C 3:03 TU'
$ws.Range("E21").Value = '2026-02-17 03:49:54'
$ws.Range("J21").Value = '1015.4 hPa'
$ws.Range("E22").Value = '2026-02-17 03:49:57'
$ws.Range("E23").Value = '2026-02-17 03:49:59'
$ws.Range("H23").Value = '65%'
$ws.Range("I23").Value = '0.4 mm'
$ws.Range("E24").Value = '2026-02-17 03:50:02'
$ws.Range("L24").Value = '43.6 km/h - 284º 3:02 TU'
$ws.Range("O24").Value = '10.0 °C'
$ws.Range("E25").Value = '2026-02-17 03:50:04'
$ws.Range("I25").Value = '0.4 mm'
$ws.Range("E26").Value = '2026-02-17 03:50:07'
$ws.Range("E27").Value = '2026-02-17 03:50:09'
$ws.Range("O27").Value = '-2.4 °C'
$ws.Range("E28").Value = '2026-02-17 03:50:12'
$ws.Range("H28").Value = '90%'
$ws.Range("J28").Value = '1015.0 hPa'
$ws.Range("O28").Value = '5.0 °C'
$ws.Range("E29").Value = '2026-02-17 03:50:15'
$ws.Range("H29").Value = '62%'
$ws.Range("L29").Value = '20.9 km/h - 350º 3:14 TU'
$ws.Range("M29").Value = '12.6 °C 3:16 TU'
$ws.Range("O29").Value = '11.5 °C'
$ws.Range("E30").Value = '2026-02-17 03:50:17'
$ws.Range("H30").Value = '46%'
$ws.Range("L30").Value = '55.4 km/h - 6º 3:21 TU'
$ws.Range("E31").Value = '2026-02-17 03:50:20'
$ws.Range("N31").Value = '9.1 °C 3:29 TU'
$ws.Range("O31").Value = '10.5 °C'
$ws.Range("E32").Value = '2026-02-17 03:50:23'
$ws.Range("N32").Value = '6.0 °C 3:23 TU'
$ws.Range("O32").Value = '6.7 °C'
$ws.Range("E33").Value = '2026-02-17 03:50:25'
$ws.Range("J33").Value = '1015.4 hPa'
$ws.Range("E34").Value = '2026-02-17 03:50:28'
$ws.Range("H34").Value = '47%'
$ws.Range("L34").Value = '56.2 km/h - 41º 3:16 TU'
$ws.Range("E35").Value = '2026-02-17 03:50:31'
$ws.Range("I35").Value = '2.0 mm'
$ws.Range("N35").Value = '4.8 °C 3:26 TU'
$ws.Range("O35").Value = '5.7 °C'
$ws.Range("E36").Value = '2026-02-17 03:50:34'
$ws.Range("H36").Value = '44%'
$ws.Range("J36").Value = '1015.2 hPa'
$ws.Range("O36").Value = '13.6 °C'
$ws.Range("E37").Value = '2026-02-17 03:50:37'
$ws.Range("H37").Value = '31%'
$ws.Range("J37").Value = '1014.7 hPa'
$ws.Range("O37").Value = '9.0 °C'
$ws.Range("E38").Value = '2026-02-17 03:50:39'
$ws.Range("O38").Value = '8.8 °C'
$ws.Range("E39").Value = '2026-02-17 03:50:42'
$ws.Range("H39").Value = '59%'
$ws.Range("I39").Value = '0.3 mm'
$ws.Range("E40").Value = '2026-02-17 03:50:45'
$ws.Range("H40").Value = '60%'
$ws.Range("J40").Value = '1017.1 hPa'
$ws.Range("O40").Value = '5.9 °C'
$ws.Range("E41").Value = '2026-02-17 03:50:47'
$ws.Range("N41").Value = '14.4 °C 3:12 TU'
$ws.Range("O41").Value = '15.6 °C'
$ws.Range("E42").Value = '2026-02-17 03:50:50'
$ws.Range("O42").Value = '13.0 °C'
$ws.Range("E43").Value = '2026-02-17 03:50:52'
$ws.Range("N43").Value = '3.6 °C 3:19 TU'
$ws.Range("O43").Value = '5.0 °C'
$ws.Range("E44").Value = '2026-02-17 03:50:55'
$ws.Range("H44").Value = '73%'
$ws.Range("E45").Value = '2026-02-17 03:50:57'
$ws.Range("H45").Value = '56%'
$ws.Range("M45").Value = '5.6 °C 3:15 TU'
$ws.Range("O45").Value = '4.4 °C'
$ws.Range("E46").Value = '2026-02-17 03:51:00'
$ws.Range("H46").Value = '56%'
$ws.Range("J46").Value = '1017.4 hPa'
$ws.Range("N46").Value = '12.7 °C 3:23 TU'
$ws.Range("O46").Value = '13.8 °C'
